$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently lists contingencies line1..line6 (rows 2-7) followed by
# extr1..extr8 (rows 8-15). Two new line contingencies (line7, line8) are
# inserted after line6, so extr1..extr8 shift down two rows (to rows 10-17).
# Copy the existing extr rows downward first (bottom-up so sources aren't
# clobbered before they're read).
for ($r = 15; $r -ge 8; $r--) {
    $destRow = $r + 2
    $ws.Cells.Item($destRow, 1).Value = $ws.Cells.Item($r, 1).Value2
    $ws.Cells.Item($destRow, 2).Value = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($destRow, 3).Value = $ws.Cells.Item($r, 3).Value2
    $ws.Cells.Item($destRow, 4).Value = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($destRow, 5).Value = $ws.Cells.Item($r, 5).Value2
}

# The two new rows (16 & 17) sit past the old used range, so they need the
# same formatting (bold/border/centered) as the rest of column A; copy it
# over from an existing formatted cell instead of building it up property by
# property (which would mint extra unused style entries).
$ws.Range("A8").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)

# Row 8: new contingency "line7"
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "line7"
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = $true

# Row 9: new contingency "line8"
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "line8"
$ws.Range("C9").Value = 16
$ws.Range("D9").Value = 9
$ws.Range("E9").Value = $true

# Renumber the index column for the shifted extr rows (A10..A17 = 8..15)
$ws.Range("A10").Value = 8
$ws.Range("A11").Value = 9
$ws.Range("A12").Value = 10
$ws.Range("A13").Value = 11
$ws.Range("A14").Value = 12
$ws.Range("A15").Value = 13
$ws.Range("A16").Value = 14
$ws.Range("A17").Value = 15

# A few of the shifted extr rows also got their "in_service" flag flipped.
$ws.Range("E13").Value = $false
$ws.Range("E15").Value = $false
$ws.Range("E16").Value = $true
